$wb = $excel.ActiveWorkbook

# --- Select Z7 on the existing "CNN 1D" sheet (was AA3:AA51) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$excel.ActiveWindow.ScrollRow = 1
$ws1.Range("Z7").Select()

# --- Add the new "CNN 2D" sheet after "CNN 1D" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "CNN 2D"

# --- Headers (row 1): "298 Packets" / "290 Packets" ---
$ws2.Range("B1").Value = "298 Packets"
$ws2.Range("F1").Value = "290 Packets"

# --- Headers (row 2): Acc / Loss / Time, twice ---
$ws2.Range("A2").Value = "Acc"
$ws2.Range("B2").Value = "Loss"
$ws2.Range("C2").Value = "Time"
$ws2.Range("E2").Value = "Acc"
$ws2.Range("F2").Value = "Loss"
$ws2.Range("G2").Value = "Time"

# --- Bold the header cells (matches style index 1 used on "CNN 1D") ---
$ws2.Range("B1").Font.Bold = $true
$ws2.Range("F1").Font.Bold = $true
$ws2.Range("A2").Font.Bold = $true
$ws2.Range("B2").Font.Bold = $true
$ws2.Range("C2").Font.Bold = $true
$ws2.Range("E2").Font.Bold = $true
$ws2.Range("F2").Font.Bold = $true
$ws2.Range("G2").Font.Bold = $true

# --- Data rows 3-51 (columns A-C; D is left blank, E-G unused) ---
$ws2.Range("A3").Value = 87.468874454498206
$ws2.Range("B3").Value = 0.27769745134200002
$ws2.Range("C3").Value = 89.2846999168396
$ws2.Range("A4").Value = 87.563735246658297
$ws2.Range("B4").Value = 0.29503988322422398
$ws2.Range("C4").Value = 85.845623970031696
$ws2.Range("A5").Value = 88.000094890594397
$ws2.Range("B5").Value = 0.26785905080514
$ws2.Range("C5").Value = 85.7653968334198
$ws2.Range("A6").Value = 87.767684459686194
$ws2.Range("B6").Value = 0.26710454454858301
$ws2.Range("C6").Value = 85.182788848876896
$ws2.Range("A7").Value = 83.306849002838106
$ws2.Range("B7").Value = 0.28680863579486499
$ws2.Range("C7").Value = 82.337363481521606
$ws2.Range("A8").Value = 87.162947654724107
$ws2.Range("B8").Value = 0.28182332678524202
$ws2.Range("C8").Value = 81.9637322425842
$ws2.Range("A9").Value = 87.698912620544405
$ws2.Range("B9").Value = 0.28215008380884199
$ws2.Range("C9").Value = 82.482324123382497
$ws2.Range("A10").Value = 83.330565690994206
$ws2.Range("B10").Value = 0.291507470804097
$ws2.Range("C10").Value = 86.969259738922105
$ws2.Range("A11").Value = 87.459385395050006
$ws2.Range("B11").Value = 0.274819421176835
$ws2.Range("C11").Value = 88.251056432723999
$ws2.Range("A12").Value = 87.826973199844304
$ws2.Range("B12").Value = 0.26558895921090497
$ws2.Range("C12").Value = 86.576238393783498
$ws2.Range("A13").Value = 87.976378202438298
$ws2.Range("B13").Value = 0.26511989357827398
$ws2.Range("C13").Value = 87.722548484802203
$ws2.Range("A14").Value = 87.537646293640094
$ws2.Range("B14").Value = 0.28136088573315199
$ws2.Range("C14").Value = 81.793133020401001
$ws2.Range("A15").Value = 87.122631072998004
$ws2.Range("B15").Value = 0.30348183921715999
$ws2.Range("C15").Value = 86.800741195678697
$ws2.Range("A16").Value = 88.555032014846802
$ws2.Range("B16").Value = 0.26193962489184303
$ws2.Range("C16").Value = 84.190393209457397
$ws2.Range("A17").Value = 83.700525760650606
$ws2.Range("B17").Value = 0.29319895227982401
$ws2.Range("C17").Value = 85.059386491775498
$ws2.Range("A18").Value = 87.964522838592501
$ws2.Range("B18").Value = 0.28085267110243001
$ws2.Range("C18").Value = 88.176538705825806
$ws2.Range("A19").Value = 87.774801254272404
$ws2.Range("B19").Value = 0.28680695181275501
$ws2.Range("C19").Value = 82.430211544036794
$ws2.Range("A20").Value = 87.732112407684298
$ws2.Range("B20").Value = 0.275884313022213
$ws2.Range("C20").Value = 85.980497121810899
$ws2.Range("A21").Value = 83.211988210678101
$ws2.Range("B21").Value = 0.29968467387745201
$ws2.Range("C21").Value = 87.220519542694007
$ws2.Range("A22").Value = 88.185071945190401
$ws2.Range("B22").Value = 0.29059914103836698
$ws2.Range("C22").Value = 86.060844182968097
$ws2.Range("A23").Value = 82.941639423370304
$ws2.Range("B23").Value = 2.6234215325207999
$ws2.Range("C23").Value = 82.944344520568805
$ws2.Range("A24").Value = 87.746340036392198
$ws2.Range("B24").Value = 0.26187385342304698
$ws2.Range("C24").Value = 83.588365554809499
$ws2.Range("A25").Value = 82.941639423370304
$ws2.Range("B25").Value = 2.6234215325207999
$ws2.Range("C25").Value = 85.223354816436697
$ws2.Range("A26").Value = 82.941639423370304
$ws2.Range("B26").Value = 2.6234215325207999
$ws2.Range("C26").Value = 81.924137830734196
$ws2.Range("A27").Value = 82.941639423370304
$ws2.Range("B27").Value = 0.32914232458174197
$ws2.Range("C27").Value = 88.039988994598303
$ws2.Range("A28").Value = 87.770056724548297
$ws2.Range("B28").Value = 0.26239407995377001
$ws2.Range("C28").Value = 86.666295766830402
$ws2.Range("A29").Value = 86.963737010955796
$ws2.Range("B29").Value = 0.28287483567024402
$ws2.Range("C29").Value = 86.143874883651705
$ws2.Range("A30").Value = 87.5661075115203
$ws2.Range("B30").Value = 0.28233318821814501
$ws2.Range("C30").Value = 82.168639659881507
$ws2.Range("A31").Value = 86.911565065383897
$ws2.Range("B31").Value = 0.29228865190134301
$ws2.Range("C31").Value = 82.568189859390202
$ws2.Range("A32").Value = 87.355041503906193
$ws2.Range("B32").Value = 0.297896636297244
$ws2.Range("C32").Value = 88.141471385955796
$ws2.Range("A33").Value = 88.367682695388794
$ws2.Range("B33").Value = 0.26378162383908099
$ws2.Range("C33").Value = 86.535898685455294
$ws2.Range("A34").Value = 88.713926076888995
$ws2.Range("B34").Value = 0.26751537651942398
$ws2.Range("C34").Value = 88.09033203125
$ws2.Range("A35").Value = 87.606424093246403
$ws2.Range("B35").Value = 0.27296707574541801
$ws2.Range("C35").Value = 87.8489665985107
$ws2.Range("A36").Value = 88.611948490142794
$ws2.Range("B36").Value = 0.27581942583646002
$ws2.Range("C36").Value = 86.214573621749807
$ws2.Range("A37").Value = 87.509191036224294
$ws2.Range("B37").Value = 0.28830567644077898
$ws2.Range("C37").Value = 82.366712570190401
$ws2.Range("A38").Value = 82.941639423370304
$ws2.Range("B38").Value = 2.6234215325207999
$ws2.Range("C38").Value = 87.390394210815401
$ws2.Range("A39").Value = 87.992978096008301
$ws2.Range("B39").Value = 0.261784661588019
$ws2.Range("C39").Value = 82.509093046188298
$ws2.Range("A40").Value = 87.402468919754
$ws2.Range("B40").Value = 0.27729379129241
$ws2.Range("C40").Value = 86.764262437820406
$ws2.Range("A41").Value = 88.014322519302297
$ws2.Range("B41").Value = 0.264681776859853
$ws2.Range("C41").Value = 85.485150337219196
$ws2.Range("A42").Value = 84.355062246322603
$ws2.Range("B42").Value = 0.31403404351810699
$ws2.Range("C42").Value = 88.892863273620605
$ws2.Range("A43").Value = 86.766904592513995
$ws2.Range("B43").Value = 0.276541047344738
$ws2.Range("C43").Value = 83.066589593887301
$ws2.Range("A44").Value = 82.941639423370304
$ws2.Range("B44").Value = 2.6234215325207999
$ws2.Range("C44").Value = 81.904806852340698
$ws2.Range("A45").Value = 82.941639423370304
$ws2.Range("B45").Value = 0.33979109252344702
$ws2.Range("C45").Value = 82.287146091461096
$ws2.Range("A46").Value = 87.682312726974402
$ws2.Range("B46").Value = 0.261655265279719
$ws2.Range("C46").Value = 82.849242687225299
$ws2.Range("A47").Value = 86.973226070403996
$ws2.Range("B47").Value = 0.30061166674435102
$ws2.Range("C47").Value = 84.153202295303302
$ws2.Range("A48").Value = 87.281525135040198
$ws2.Range("B48").Value = 0.29289897332598203
$ws2.Range("C48").Value = 88.161750078201294
$ws2.Range("A49").Value = 87.260180711746202
$ws2.Range("B49").Value = 0.29638089413215002
$ws2.Range("C49").Value = 86.880231142043996
$ws2.Range("A50").Value = 82.941639423370304
$ws2.Range("B50").Value = 2.6234215325207999
$ws2.Range("C50").Value = 87.098807334899902
$ws2.Range("A51").Value = 87.736856937408405
$ws2.Range("B51").Value = 0.271246363287586
$ws2.Range("C51").Value = 82.711402893066406

# --- Set the selection on the new sheet and make it the active tab ---
$ws2.Range("N10").Select()
$ws2.Activate()

Write-Host "done"
